$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B2").Value = "701号直流"
$ws.Range("C2").Value = 45927.457337962966
$ws.Range("D2").Value = 45932.31215277778
$ws.Range("E2").Value = 116.51555555546656

$ws.Range("A3").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B3").Value = "502号直流"
$ws.Range("C3").Value = 45930.23804398148
$ws.Range("D3").Value = 45932.31215277778
$ws.Range("E3").Value = 49.77861111110542

$ws.Range("A4").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B4").Value = "001A号直流"
$ws.Range("C4").Value = 45930.24489583333
$ws.Range("D4").Value = 45932.31215277778
$ws.Range("E4").Value = 49.614166666637175

$ws.Range("A5").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B5").Value = "112号直流"
$ws.Range("C5").Value = 45930.517060185186
$ws.Range("D5").Value = 45932.31215277778
$ws.Range("E5").Value = 43.082222222175915

$ws.Range("A6").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B6").Value = "904号直流"
$ws.Range("C6").Value = 45930.554375
$ws.Range("D6").Value = 45932.31215277778
$ws.Range("E6").Value = 42.18666666664649

$ws.Range("A7").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B7").Value = "101号直流"
$ws.Range("C7").Value = 45930.56182870371
$ws.Range("D7").Value = 45932.31215277778
$ws.Range("E7").Value = 42.007777777675074

$ws.Range("A8").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B8").Value = "111号直流"
$ws.Range("C8").Value = 45930.61851851852
$ws.Range("D8").Value = 45932.31215277778
$ws.Range("E8").Value = 40.64722222223645

$ws.Range("A9").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B9").Value = "103号直流"
$ws.Range("C9").Value = 45930.64413194444
$ws.Range("D9").Value = 45932.31215277778
$ws.Range("E9").Value = 40.03250000003027

$ws.Range("A10").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B10").Value = "306号直流"
$ws.Range("C10").Value = 45930.69081018519
$ws.Range("D10").Value = 45932.31215277778
$ws.Range("E10").Value = 38.912222222134005

$ws.Range("A11").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B11").Value = "105号直流"
$ws.Range("C11").Value = 45931.039722222224
$ws.Range("D11").Value = 45932.31215277778
$ws.Range("E11").Value = 30.53833333327202

$ws.Range("A12").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B12").Value = "602号直流"
$ws.Range("C12").Value = 45931.22626157408
$ws.Range("D12").Value = 45932.31215277778
$ws.Range("E12").Value = 26.061388888803776

$ws.Range("A13").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B13").Value = "604号直流"
$ws.Range("C13").Value = 45931.22881944444
$ws.Range("D13").Value = 45932.31215277778
$ws.Range("E13").Value = 26.000000000058208

$ws.Range("A14").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B14").Value = "705号直流"
$ws.Range("C14").Value = 45931.23175925926
$ws.Range("D14").Value = 45932.31215277778
$ws.Range("E14").Value = 25.929444444424007

$ws.Range("A15").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B15").Value = "206号直流"
$ws.Range("C15").Value = 45931.240578703706
$ws.Range("D15").Value = 45932.31215277778
$ws.Range("E15").Value = 25.71777777769603

$ws.Range("A16").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B16").Value = "503号直流"
$ws.Range("C16").Value = 45931.260625
$ws.Range("D16").Value = 45932.31215277778
$ws.Range("E16").Value = 25.23666666657664

$ws.Range("A17").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B17").Value = "107号直流"
$ws.Range("C17").Value = 45931.411145833335
$ws.Range("D17").Value = 45932.31215277778
$ws.Range("E17").Value = 21.62416666658828

$ws.Range("A18").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B18").Value = "005A号直流"
$ws.Range("C18").Value = 45931.51777777778
$ws.Range("D18").Value = 45932.31215277778
$ws.Range("E18").Value = 19.06499999994412

$ws.Range("A19").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B19").Value = "108号直流"
$ws.Range("C19").Value = 45931.51920138889
$ws.Range("D19").Value = 45932.31215277778
$ws.Range("E19").Value = 19.030833333323244

$ws.Range("A20").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B20").Value = "006B号直流"
$ws.Range("C20").Value = 45931.53244212963
$ws.Range("D20").Value = 45932.31215277778
$ws.Range("E20").Value = 18.713055555475876

$ws.Range("A21").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B21").Value = "002A号直流"
$ws.Range("C21").Value = 45931.55152777778
$ws.Range("D21").Value = 45932.31215277778
$ws.Range("E21").Value = 18.25499999988824

$ws.Range("A22").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B22").Value = "403号直流"
$ws.Range("C22").Value = 45931.55615740741
$ws.Range("D22").Value = 45932.31215277778
$ws.Range("E22").Value = 18.143888888822403

$ws.Range("A23").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B23").Value = "201号直流"
$ws.Range("C23").Value = 45931.57554398148
$ws.Range("D23").Value = 45932.31215277778
$ws.Range("E23").Value = 17.678611111070495

$ws.Range("A24").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B24").Value = "210号直流"
$ws.Range("C24").Value = 45931.584861111114
$ws.Range("D24").Value = 45932.31215277778
$ws.Range("E24").Value = 17.454999999899883

$ws.Range("A25").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B25").Value = "203号直流"
$ws.Range("C25").Value = 45931.585393518515
$ws.Range("D25").Value = 45932.31215277778
$ws.Range("E25").Value = 17.44222222227836

$ws.Range("A26").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B26").Value = "702号直流"
$ws.Range("C26").Value = 45931.58704861111
$ws.Range("D26").Value = 45932.31215277778
$ws.Range("E26").Value = 17.40250000002561

$ws.Range("A27").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B27").Value = "401号直流"
$ws.Range("C27").Value = 45931.59206018518
$ws.Range("D27").Value = 45932.31215277778
$ws.Range("E27").Value = 17.282222222245764

$ws.Range("A28").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B28").Value = "A01号直流"
$ws.Range("C28").Value = 45931.60625
$ws.Range("D28").Value = 45932.31215277778
$ws.Range("E28").Value = 16.941666666709352

$ws.Range("A29").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B29").Value = "B02号直流"
$ws.Range("C29").Value = 45931.60697916667
$ws.Range("D29").Value = 45932.31215277778
$ws.Range("E29").Value = 16.924166666634846

$ws.Range("A30").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B30").Value = "304号直流"
$ws.Range("C30").Value = 45931.6202662037
$ws.Range("D30").Value = 45932.31215277778
$ws.Range("E30").Value = 16.605277777765878

$ws.Range("A31").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B31").Value = "301号直流"
$ws.Range("C31").Value = 45931.62170138889
$ws.Range("D31").Value = 45932.31215277778
$ws.Range("E31").Value = 16.57083333330229

$ws.Range("A32").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B32").Value = "401号直流"
$ws.Range("C32").Value = 45931.64065972222
$ws.Range("D32").Value = 45932.31215277778
$ws.Range("E32").Value = 16.11583333328599

$ws.Range("A33").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B33").Value = "505号直流"
$ws.Range("C33").Value = 45931.6472337963
$ws.Range("D33").Value = 45932.31215277778
$ws.Range("E33").Value = 15.95805555547122

$ws.Range("A34").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B34").Value = "102号直流"
$ws.Range("C34").Value = 45931.6478125
$ws.Range("D34").Value = 45932.31215277778
$ws.Range("E34").Value = 15.944166666653473

$ws.Range("A35").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B35").Value = "007A号直流"
$ws.Range("C35").Value = 45931.666666666664
$ws.Range("D35").Value = 45932.31215277778
$ws.Range("E35").Value = 15.49166666669771

$ws.Range("A36").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B36").Value = "805号直流"
$ws.Range("C36").Value = 45931.766435185185
$ws.Range("D36").Value = 45932.31215277778
$ws.Range("E36").Value = 13.097222222189885

$ws.Range("A37").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B37").Value = "106号直流"
$ws.Range("C37").Value = 45931.789456018516
$ws.Range("D37").Value = 45932.31215277778
$ws.Range("E37").Value = 12.544722222257406

$ws.Range("A38").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B38").Value = "101号直流"
$ws.Range("C38").Value = 45931.80849537037
$ws.Range("D38").Value = 45932.31215277778
$ws.Range("E38").Value = 12.087777777691372

$ws.Range("K10").Select()